$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 66, shifting existing rows 66-69 down to 67-70.
$ws.Rows.Item(66).Insert()

# Populate the newly inserted row 66 with the new weekly record.
$ws.Cells.Item(66, 1).Value = 4
$ws.Cells.Item(66, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(66, 3).Value = "Los Lagos"
$ws.Cells.Item(66, 4).Value = 44995
$ws.Cells.Item(66, 5).Value = 10
$ws.Cells.Item(66, 6).Value = 100112030
$ws.Cells.Item(66, 7).Value = "Poroto granado"
$ws.Cells.Item(66, 8).Value = "Sin especificar"
$ws.Cells.Item(66, 9).Value = "Primera"
$ws.Cells.Item(66, 10).Value = 70
$ws.Cells.Item(66, 11).Value = 36000
$ws.Cells.Item(66, 12).Value = 36000
$ws.Cells.Item(66, 13).Value = 36000
$ws.Cells.Item(66, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(66, 15).Value = "Región del Maule"
$ws.Cells.Item(66, 16).Value = 1440
$ws.Cells.Item(66, 17).Value = 25
$ws.Cells.Item(66, 18).Value = "Hortaliza"
